# Append "/fgfgf" (Arial) after the existing Hindi transliteration text
# (Lohit Hindi) in cell C2, as two distinct rich-text runs, then move the
# active selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("C2")

# Original text was: "एजुकेशनल टेक्नोलॉजी" (19 characters, rendered with the
# "Lohit Hindi" font via the cell's style). New text appends "/fgfgf" in
# Arial, so the cell keeps two runs with different fonts.
$original = "एजुकेशनल टेक्नोलॉजी"
$suffix = "/fgfgf"
$cell.Value = $original + $suffix

$originalLen = $original.Length
$suffixLen = $suffix.Length

# Run 1: the pre-existing Hindi text, kept in "Lohit Hindi".
$run1 = $cell.Characters(1, $originalLen)
$run1.Font.Name = "Lohit Hindi"
$run1.Font.Size = 10
$run1.Font.ColorIndex = -4105

# Run 2: the newly appended text, in Arial.
$run2 = $cell.Characters($originalLen + 1, $suffixLen)
$run2.Font.Name = "Arial"
$run2.Font.Size = 10
$run2.Font.ColorIndex = -4105

# Move the active selection to the edited cell.
[void]$ws.Range("C2").Select()
